$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows for the new weekly report (week of 2021-11-16, serial 44516),
# pushing the existing rows 89-97 down to 95-103.
$ws.Rows("89:94").Insert()

$data = @(
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44516, 13, 300000000, "Espárragos", "Sin especificar", "Banquete", 270, 1400, 1500, 1456, "`$/kilo", "Provincia de Linares", 1456, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44516, 13, 300000000, "Espárragos", "Sin especificar", "Banquete", 580, 1200, 1300, 1260, "`$/kilo", "Región Metropolitana", 1260, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44516, 13, 300000000, "Espárragos", "Sin especificar", "Primera", 300, 1200, 1300, 1233, "`$/kilo", "Provincia de Linares", 1233, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44516, 13, 300000000, "Espárragos", "Sin especificar", "Primera", 720, 1000, 1100, 1076, "`$/kilo", "Región Metropolitana", 1076, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44516, 13, 300000000, "Espárragos", "Sin especificar", "Segunda", 140, 1000, 1100, 1036, "`$/kilo", "Provincia de Linares", 1036, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44516, 13, 300000000, "Espárragos", "Sin especificar", "Segunda", 300, 800, 900, 867, "`$/kilo", "Región Metropolitana", 867, 1, "Hortaliza")
)

$row = 89
foreach ($rowData in $data) {
    $col = 1
    foreach ($val in $rowData) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}

# Ensure the date column keeps the date number format used elsewhere in column D.
$ws.Range("D89:D94").NumberFormat = $ws.Range("D88").NumberFormat
